$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: response changes to FALSE, confidence becomes 0.8
$ws.Range("B2").Value = "'FALSE"
$ws.Range("C2").Value = 0.8

# Row 3: confidence changes to 0.8
$ws.Range("C3").Value = 0.8

# Row 4: confidence changes to 1
$ws.Range("C4").Value = 1

# Row 5: response changes to FALSE, confidence becomes 0.8
$ws.Range("B5").Value = "'FALSE"
$ws.Range("C5").Value = 0.8

# Row 8: response changes to FALSE, confidence becomes 0.8
$ws.Range("B8").Value = "'FALSE"
$ws.Range("C8").Value = 0.8

# Row 10: confidence changes to 0.9
$ws.Range("C10").Value = 0.9
